$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.533.28"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "2.918.31"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "498.47"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.93"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.427"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.02"
$ws.Range("E9").Value = "  -3.69%  "
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.363"
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D12").Value = "3.491.63"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.99"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000160"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "55.635.44"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "2.973.55"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.92"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.83"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.79"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.19"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.994"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.487"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.73"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "3.119.93"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.158"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").Value = "0.0₃0873"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.45"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.77"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.98"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.30"
$ws.Range("E34").Value = "  -2.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.53"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.72"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.24"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.58"
$ws.Range("E38").Value = "  +2.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0653"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "2.972.09"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.52"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.642"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").Value = "2.151.65"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.35"
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.85"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.921"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0236"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.45"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0844"
$ws.Range("E51").Value = "  -3.27%  "
